# Updates the cryptos price/volume table to the latest scraped values.
# (Generated for commit: "Updated cryptos list on Sat Sep 23 14:30:01 UTC 2023 with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell {
    param(
        [string]$Ref,
        [string]$Value
    )
    $ws.Range($Ref).Value = $Value
}

# Row 2 - Bitcoin
Set-Cell "D2" "26.668.59"
Set-Cell "E2" "  +0.02%  "

# Row 3 - Ethereum
Set-Cell "D3" "1.596.39"
Set-Cell "E3" "  -0.06%  "

# Row 4 - TetherUSD
Set-Cell "E4" "  +0.28%  "

# Row 5 - BNB
Set-Cell "D5" "211.45"
Set-Cell "E5" "  +0.34%  "

# Row 6 - XRP
Set-Cell "E6" "  +0.02%  "

# Row 7 - USDC
Set-Cell "E7" "  +0.22%  "

# Row 8 - Dogecoin
Set-Cell "E8" "  +0.19%  "

# Row 9 - Cardano
Set-Cell "D9" "0.247"
Set-Cell "E9" "  +0.42%  "

# Row 10 - Solana
Set-Cell "D10" "19.50"
Set-Cell "E10" "  -0.55%  "

# Row 11 - TRON
Set-Cell "E11" "  -0.33%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-Cell "D12" "1.820.60"
Set-Cell "E12" "  -0.04%  "

# Row 13 - WrappedEther
Set-Cell "D13" "1.607.84"
Set-Cell "E13" "  +1.14%  "

# Row 14 - Polkadot
Set-Cell "D14" "4.03"
Set-Cell "E14" "  +0.27%  "

# Row 15 - Polygon
Set-Cell "E15" "  +0.64%  "

# Row 16 - Litecoin
Set-Cell "D16" "65.02"
Set-Cell "E16" "  +0.44%  "

# Row 17 - WrappedBTC
Set-Cell "D17" "26.641.28"
Set-Cell "E17" "  -0.01%  "

# Row 18 - ShibaInu
Set-Cell "D18" "0.0`u{2083}0745"
Set-Cell "E18" "  +2.31%  "

# Rows 19/20 - Dai and BitcoinCash swap positions (BitcoinCash now ranked above Dai)
Set-Cell "B19" "BitcoinCash"
Set-Cell "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-Cell "D19" "209.21"
Set-Cell "E19" "  +0.40%  "

Set-Cell "B20" "Dai"
Set-Cell "C20" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-Cell "D20" "1.00"
Set-Cell "E20" "  +0.24%  "

# Row 21 - Chainlink
Set-Cell "D21" "7.06"
Set-Cell "E21" "  +4.41%  "

# Row 22 - Uniswap
Set-Cell "E22" "  +1.04%  "

# Row 23 - Toncoin
Set-Cell "D23" "2.33"
Set-Cell "E23" "  +0.98%  "

# Row 24 - Avalanche
Set-Cell "D24" "8.98"
Set-Cell "E24" "  +0.90%  "

# Row 25 - Monero
Set-Cell "D25" "143.09"
Set-Cell "E25" "  -1.61%  "

# Row 26 - BinanceUSD
Set-Cell "E26" "  +0.15%  "

# Row 27 - Cosmos
Set-Cell "E27" "  -1.40%  "

# Row 28 - Stellar
Set-Cell "E28" "  -1.07%  "

# Row 29 - EthereumClassic
Set-Cell "D29" "15.33"
Set-Cell "E29" "  +0.52%  "

# Row 30 - Hedera
Set-Cell "D30" "0.0515"
Set-Cell "E30" "  +1.97%  "

# Row 31 - PancakeSwap
Set-Cell "E31" "  +0.20%  "

# Row 32 - Filecoin
Set-Cell "D32" "3.24"
Set-Cell "E32" "  +0.06%  "

# Row 33 - InternetComputer(DFINITY)
Set-Cell "D33" "2.95"

# Row 34 - Maker
Set-Cell "D34" "1.288.37"
Set-Cell "E34" "  -0.27%  "

# Row 35 - ImmutableX
Set-Cell "D35" "0.619"
Set-Cell "E35" "  -5.60%  "

# Row 36 - HuobiToken
Set-Cell "E36" "  +0.51%  "

# Row 37 - LidoDAOToken
Set-Cell "D37" "1.49"

# Row 38 - VeChain
Set-Cell "E38" "  -0.56%  "

# Row 39 - ARBITRUM
Set-Cell "D39" "0.829"
Set-Cell "E39" "  -1.96%  "

# Row 40 - WEMIXToken
Set-Cell "E40" "  +17.06%  "

# Row 41 - FraxShare
Set-Cell "D41" "5.44"
Set-Cell "E41" "  +0.95%  "

# Row 42 - MXToken
Set-Cell "E42" "  -0.31%  "

# Row 43 - TrustWalletToken
Set-Cell "D43" "0.781"
Set-Cell "E43" "  -0.56%  "

# Row 44 - Aave
Set-Cell "D44" "63.24"
Set-Cell "E44" "  -0.73%  "

# Row 45 - RocketPoolETH
Set-Cell "D45" "1.732.48"
Set-Cell "E45" "  -0.10%  "

# Row 46 - Quant
Set-Cell "D46" "91.13"

# Row 47 - RenderToken
Set-Cell "E47" "  -2.68%  "

# Row 48 - Algorand
Set-Cell "E48" "  +0.90%  "

# Row 50 - USDD
Set-Cell "E50" "  +0.28%  "

# Row 51 - EnergySwap
Set-Cell "E51" "  -1.07%  "
